# Reclassify parrotfishes previously labelled "Scaridae" as "Labridae"
# (Scaridae is now considered a subfamily within Labridae).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$count = 0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $c = $ws.Cells.Item($r, 1)
    $v = $c.Value()
    if ($v -eq "Scaridae") {
        $c.Value = "Labridae"
        $count++
    }
}

Write-Host "Replaced $count cells in column A (Family) from 'Scaridae' to 'Labridae'"
